$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: remove a stale <w:lastRenderedPageBreak/> marker from the run that
# starts at a paragraph's first character, WITHOUT letting the engine merge
# that run into any of its identically-formatted neighbours (which would
# happen if we simply reset the whole paragraph's text).
#
# Trick: temporarily flip Bold on every *other* run in the paragraph so it
# no longer has identical formatting to the run we are about to touch: this
# keeps that run (and all later runs) as separate, untouched runs. Then we
# force a tiny "delete char / retype char" edit on the first run only (that
# is enough for the host to drop the cached lastRenderedPageBreak flag when
# it rebuilds that run), and finally restore the Bold flag we toggled away,
# one run at a time so nothing gets re-merged either.
# ---------------------------------------------------------------------------
function Remove-LastRenderedPageBreak($paraIndex, $boundaries) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $start = $full.Start
    $firstChar = $full.Text.Substring(0, 1)

    # boundaries is a list of character offsets (relative to paragraph start)
    # delimiting run 2, run 3, ... (run 1 is [0, boundaries[0]])
    for ($i = 0; $i -lt $boundaries.Count - 1; $i++) {
        $rs = $start + $boundaries[$i]
        $re = $start + $boundaries[$i + 1]
        $rr = $d.Range($rs, $re)
        $rr.Bold = 0
    }

    $r1 = $d.Range($start, $start + 1)
    $r1.Text = "*"
    $r1b = $d.Range($start, $start + 1)
    $r1b.Text = $firstChar

    for ($i = 0; $i -lt $boundaries.Count - 1; $i++) {
        $rs = $start + $boundaries[$i]
        $re = $start + $boundaries[$i + 1]
        $rr = $d.Range($rs, $re)
        $rr.Bold = 1
    }
}

# "Table 3: Performance of quadratic probing in various load factors"
# runs: "Table " | "3" | ": Performance of " | "quadratic" | " probing in various load factors"
Remove-LastRenderedPageBreak 88 @(6, 7, 24, 33, 65)

# "Table 5: Performance of various collision resolution methods in load factor 0.4"
# single run -> no neighbours to protect
Remove-LastRenderedPageBreak 193 @()

# "Table 7: Performance of various collision resolution methods in load factor 0.6"
# runs: "Table " | "7" | ": Performance of various collision resolution methods in load factor 0." | "6"
Remove-LastRenderedPageBreak 270 @(6, 7, 78, 79)

# "Table 9: Performance of various collision resolution methods in load factor 0.8"
# runs: "Table " | "9" | ": Performance of various collision resolution methods in load factor 0." | "8"
Remove-LastRenderedPageBreak 347 @(6, 7, 78, 79)

# ---------------------------------------------------------------------------
# Merge the split "0.000108238m" / "s" runs back into a single run reading
# "0.000108238ms" (the table cell that uses 12pt Times New Roman - there is
# a second, unrelated "0.000108238ms" elsewhere in 11pt that must stay put).
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $rng = $p.Range
    if ($rng.Text -like "0.000108238ms*" -and $rng.Font.Size -eq 12) {
        $start = $rng.Start
        $r1 = $d.Range($start, $start + 1)
        $r1.Text = "*"
        $r1b = $d.Range($start, $start + 1)
        $r1b.Text = "0"
        break
    }
}

# ---------------------------------------------------------------------------
# Insert a new centred, bold "N = 1000003" paragraph before the very first
# paragraph of the document, copying its paragraph/run formatting.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1).Range
$firstPara.InsertParagraphBefore()
$newPara = $d.Paragraphs(1).Range
$newPara.Text = "N = 1000003"
$newPara.Font.Color = 0

Write-Host "done"
